# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, matching the commit's refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.352.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.715.88"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5302"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06689"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2654"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.94"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07695"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.489"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.952.96"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.710.24"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8191"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.77"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.377.70"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.70"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.650"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.031"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.709"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1208"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.251"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05381"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.297"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.483"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.399"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.637"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9528"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  -1.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5898"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.158.49"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +10.66%  "

$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.833"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.16%  "

$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8404"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.859.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("E46").Value = "  +4.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.81"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4570"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.167"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05199"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.97%  "
